$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.724.16"
$ws.Range("E2").Value = "  -1.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.546.13"
$ws.Range("E3").Value = "  -1.76%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.74"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("E6").Value = "  -1.93%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.43"
$ws.Range("E8").Value = "  -3.82%  "

$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("E10").Value = "  -1.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("E11").Value = "  -2.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.766.52"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.553.51"
$ws.Range("E13").Value = "  -1.47%  "

$ws.Range("E14").Value = "  -2.62%  "

$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.711.51"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.25"
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.95"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -2.19%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.94"
$ws.Range("E23").Value = "  -5.30%  "

$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.88"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("E26").Value = "  -3.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.87"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  -1.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0460"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("E32").Value = "  +0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.336.69"
$ws.Range("E33").Value = "  -4.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("E35").Value = "  -3.18%  "

$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.935"
$ws.Range("E37").Value = "  -1.22%  "

$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("E39").Value = "  +1.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  +6.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.800"
$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.65"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("E45").Value = "  -4.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.680.06"
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.24"
$ws.Range("E47").Value = "  -3.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.98"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0975"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0952"
$ws.Range("E51").Value = "  -0.07%  "
